# Edit: add a test case to the "Concepts" sheet of the ontology test workbook
# that checks concepts defined twice in the same excel sheet are reported,
# and fix the wording of the neighbouring "already exists" test description.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# 1) Fix a typo in the existing "already exists" test description (row 23, col J):
#    "...already exists (in and imported ontology)" -> "...already exists in and imported ontology)"
$ws.Range("J23").Value2 = "Test adding concept that already exists in and imported ontology)"

# 2) Append a new row (24) duplicating the "Pattern" concept defined earlier on
#    row 7, so that the same concept is defined twice in the sheet - this is the
#    new test case for duplicate-concept detection.
$ws.Range("A24").Value2 = $ws.Range("A7").Value2
$ws.Range("C24").Value2 = $ws.Range("C7").Value2
$ws.Range("C24").WrapText = $true
$ws.Range("D24").WrapText = $true
$ws.Range("G24").Value2 = $ws.Range("G7").Value2
$ws.Range("H24").WrapText = $true
$ws.Range("J24").Value2 = "Test defining same concept twice in the same excel sheet"
$ws.Rows.Item(24).RowHeight = 30

# Update the active selection to reflect where the author ended up after typing
# the new row (mirrors the recorded end-state in the source file).
$ws.Range("L26").Select()
